# Weekly fruit/vegetable price update.
# Insert a new weekly record before the current last-but-one row (162),
# pushing the existing rows 162-163 down to 163-164, and populate the
# newly inserted row 162 with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(162).Insert()

$ws.Cells.Item(162, 1).Value = 4
$ws.Cells.Item(162, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(162, 3).Value = "Los Lagos"
$ws.Cells.Item(162, 4).Value = 44656
$ws.Cells.Item(162, 5).Value = 10
$ws.Cells.Item(162, 6).Value = 100112009
$ws.Cells.Item(162, 7).Value = "Acelga"
$ws.Cells.Item(162, 8).Value = "Sin especificar"
$ws.Cells.Item(162, 9).Value = "Primera"
$ws.Cells.Item(162, 10).Value = 60
$ws.Cells.Item(162, 11).Value = 10000
$ws.Cells.Item(162, 12).Value = 10000
$ws.Cells.Item(162, 13).Value = 10000
$ws.Cells.Item(162, 14).Value = "`$/docena de atados (12 kilos)"
$ws.Cells.Item(162, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(162, 16).Value = 833
$ws.Cells.Item(162, 17).Value = 12
$ws.Cells.Item(162, 18).Value = "Hortaliza"
